{"js": "// Applies the \"docs/\u05d3\u05e3 \u05d4\u05e2\u05e8\u05d5\u05ea.docx\" content edits described by the commit.\n// Strategy: locate each affected sentence by its exact (pre-edit) text via\n// Body.search(), then replace its contents with the revised text using\n// Range.insertText(..., Word.InsertLocation.replace). This is robust to the\n// exact paragraph/run layout and does not depend on paragraph indices.\n// The one structural change (the deleted-note sentence) also relocates the\n// hidden \"_GoBack\" bookmark so it again sits right after the word \"\u05d0\u05e9\u05e4\u05d4\".\n\nasync function replaceSentence(context, oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"...\u05d5\u05d1\u05e6\u05d3 \u05e9\u05de\u05d0\u05dc \u05d4\u05ea\u05d0\u05e8\u05d9\u05da \u05e9\u05d4\u05d5\u05db\u05e0\u05e1\u05d4 \u05d1\u05d5.\" -> \"...\u05d5\u05d1\u05e6\u05d3 \u05e9\u05de\u05d0\u05dc \u05ea\u05d0\u05e8\u05d9\u05da \u05d4\u05d5\u05e1\u05e4\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4.\"\nawait replaceSentence(\n  context,\n  \"\u05d1\u05e6\u05d3 \u05d9\u05de\u05d9\u05df \u05de\u05d5\u05e4\u05d9\u05e2\u05d4 \u05db\u05d5\u05ea\u05e8\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4, \u05d5\u05d1\u05e6\u05d3 \u05e9\u05de\u05d0\u05dc \u05d4\u05ea\u05d0\u05e8\u05d9\u05da \u05e9\u05d4\u05d5\u05db\u05e0\u05e1\u05d4 \u05d1\u05d5.\",\n  \"\u05d1\u05e6\u05d3 \u05d9\u05de\u05d9\u05df \u05de\u05d5\u05e4\u05d9\u05e2\u05d4 \u05db\u05d5\u05ea\u05e8\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4, \u05d5\u05d1\u05e6\u05d3 \u05e9\u05de\u05d0\u05dc \u05ea\u05d0\u05e8\u05d9\u05da \u05d4\u05d5\u05e1\u05e4\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4.\"\n);\n\n// 2) Clicking a note -> describe what the child sees.\nawait replaceSentence(\n  context,\n  \"\u05d0\u05dd \u05dc\u05d5\u05d7\u05e6\u05d9\u05dd \u05e2\u05dc \u05d4\u05e2\u05e8\u05d4 \u05de\u05e1\u05d5\u05d9\u05de\u05ea \u05d4\u05d9\u05dc\u05d3 \u05d9\u05db\u05d5\u05dc \u05dc\u05e8\u05d0\u05d5\u05ea \u05d0\u05ea \u05d4\u05ea\u05d5\u05db\u05df, \u05e1\u05e4\u05e8 \u05d0\u05dd \u05d4\u05d5\u05e1\u05d9\u05e3 \u05e1\u05e4\u05e8 \u05de\u05e1\u05d5\u05d9\u05d9\u05dd \u05d5\u05d2\u05dd \u05e6\u05d9\u05d5\u05e8. \u05dc\u05de\u05e9\u05dc:\",\n  \"\u05d0\u05dd \u05dc\u05d5\u05d7\u05e6\u05d9\u05dd \u05e2\u05dc \u05d4\u05e2\u05e8\u05d4 \u05de\u05e1\u05d5\u05d9\u05de\u05ea \u05d4\u05d9\u05dc\u05d3 \u05d9\u05db\u05d5\u05dc \u05dc\u05e8\u05d0\u05d5\u05ea \u05d0\u05ea \u05d4\u05ea\u05d5\u05db\u05df, \u05e9\u05dd \u05d4\u05e1\u05e4\u05e8 \u05d0\u05dd \u05d4\u05d5\u05e1\u05d9\u05e3 \u05e1\u05e4\u05e8 \u05de\u05e1\u05d5\u05d9\u05d9\u05dd \u05d5\u05d2\u05dd \u05e6\u05d9\u05d5\u05e8 \u05d0\u05dd \u05d4\u05e2\u05dc\u05d4 \u05e6\u05d9\u05d5\u05e8. \u05dc\u05de\u05e9\u05dc:\"\n);\n\n// 3) \"\u05d4\u05d5\u05e1\u05e3 \u05d4\u05e2\u05e8\u05d4\" button -> opens blank form wording.\nawait replaceSentence(\n  context,\n  \"\u05d0\u05dd \u05d4\u05d9\u05dc\u05d3 \u05dc\u05d5\u05d7\u05e5 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \\\"\u05d4\u05d5\u05e1\u05e3 \u05d4\u05e2\u05e8\u05d4\\\", \u05d4\u05d5\u05d0 \u05de\u05e7\u05d1\u05dc \u05e4\u05d5\u05e8\u05dd \u05db\u05d3\u05d9 \u05dc\u05d4\u05db\u05e0\u05d9\u05e1 \u05e4\u05e8\u05d8\u05d9 \u05d4\u05d4\u05e2\u05e8\u05d4 \u05d4\u05d7\u05d3\u05e9\u05d4:\",\n  \"\u05d0\u05dd \u05d4\u05d9\u05dc\u05d3 \u05dc\u05d5\u05d7\u05e5 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \\\"\u05d4\u05d5\u05e1\u05e3 \u05d4\u05e2\u05e8\u05d4\\\", \u05e0\u05e4\u05ea\u05d7 \u05d8\u05d5\u05e4\u05e1 \u05e8\u05d9\u05e7 \u05dc\u05d4\u05d5\u05e1\u05e4\u05ea \u05d4\u05e2\u05e8\u05d4 \u05d4\u05d7\u05d3\u05e9\u05d4:\"\n);\n\n// 4) Book list bullet wording.\nawait replaceSentence(\n  context,\n  \"\u05e1\u05e4\u05e8- \u05d1\u05d7\u05d9\u05e8\u05ea \u05e1\u05e4\u05e8 \u05de\u05e1\u05d5\u05d9\u05dd \u05de\u05e8\u05e9\u05d9\u05de\u05ea \u05d4\u05e1\u05e4\u05e8\u05d9 \u05e9\u05dc\u05d5 \u05d0\u05d5 \u05d1\u05dc\u05d9 \u05e1\u05e4\u05e8.\",\n  \"\u05e1\u05e4\u05e8- \u05d1\u05d7\u05d9\u05e8\u05ea \u05e1\u05e4\u05e8 \u05de\u05e1\u05d5\u05d9\u05dd \u05de\u05e8\u05e9\u05d9\u05de\u05ea \u05d4\u05e1\u05e4\u05e8\u05d9\u05dd \u05e9\u05dc \u05d4\u05d9\u05dc\u05d3 \u05d0\u05d5 \u05d1\u05dc\u05d9 \u05e1\u05e4\u05e8.\"\n);\n\n// 5) Pencil/update button wording.\nawait replaceSentence(\n  context,\n  \"\u05d1\u05dc\u05d7\u05d9\u05e6\u05d4 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \u05e2\u05e4\u05e8\u05d5\u05df, \u05d4\u05d9\u05dc\u05d3 \u05de\u05e7\u05d1\u05dc \u05e4\u05d5\u05e8\u05dd \u05db\u05d3\u05d9 \u05dc\u05e2\u05d3\u05db\u05df \u05d0\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4:\",\n  \"\u05d1\u05dc\u05d7\u05d9\u05e6\u05d4 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \u05e2\u05e4\u05e8\u05d5\u05df, \u05d4\u05d9\u05dc\u05d3 \u05de\u05e7\u05d1\u05dc \u05d8\u05d5\u05e4\u05e1 \u05e2\u05dc \u05de\u05e0\u05ea \u05dc\u05e2\u05d3\u05db\u05df \u05d0\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4:\"\n);\n\n// 6) Required-fields note.\nawait replaceSentence(\n  context,\n  \"\u05d7\u05d9\u05d9\u05d1 \u05e9\u05d4\u05db\u05d5\u05ea\u05e8\u05ea \u05d5\u05d4\u05ea\u05d5\u05db\u05df \u05dc\u05d0 \u05d9\u05d4\u05d9\u05d5 \u05e8\u05d9\u05e7\u05d9\u05dd.\",\n  \"\u05e9\u05d3\u05d4 \u05d4\u05db\u05d5\u05ea\u05e8\u05ea \u05d5\u05e9\u05d3\u05d4 \u05d4\u05ea\u05d5\u05db\u05df \u05dc\u05d0 \u05d9\u05d4\u05d9\u05d5 \u05e8\u05d9\u05e7\u05d9\u05dd.\"\n);\n\n// 7) Delete-note sentence, plus relocating the hidden \"_GoBack\" bookmark so\n// it ends up right after the new word \"\u05d0\u05e9\u05e4\u05d4\" (matching its position in the\n// edited document instead of the trailing empty paragraph).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nawait replaceSentence(\n  context,\n  \"\u05d1\u05dc\u05d7\u05d9\u05e6\u05d4 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \u05e4\u05d7 \u05d4\u05d4\u05e2\u05e8\u05d4 \u05e0\u05de\u05d7\u05e7\u05ea.\",\n  \"\u05d1\u05dc\u05d7\u05d9\u05e6\u05d4 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \u05d0\u05e9\u05e4\u05d4 \u05d4\u05d4\u05e2\u05e8\u05d4 \u05ea\u05d9\u05de\u05d7\u05e7.\"\n);\n\nconst trashWord = context.document.body.search(\"\u05d0\u05e9\u05e4\u05d4\", { matchCase: true });\ntrashWord.load(\"items\");\nawait context.sync();\nif (trashWord.items.length === 0) {\n  throw new Error(\"Could not find the word \u05d0\u05e9\u05e4\u05d4 after replacement\");\n}\ntrashWord.items[0].getRange(Word.RangeLocation.end).insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Applies the \"docs/\u05d3\u05e3 \u05d4\u05e2\u05e8\u05d5\u05ea.docx\" content edits described by the commit.\n# Strategy: use Find.Execute (search only, no ReplaceWith) against the full\n# document Range to locate each affected sentence by its exact pre-edit\n# text, then assign the revised text directly to Range.Text. (Using\n# Find.Execute's own ReplaceWith parameter triggers Word's smart-quotes\n# autocorrect on the replacement text, which would incorrectly curl the\n# straight quotation marks in this document - assigning Range.Text directly\n# avoids that.) This is robust to the exact paragraph/run layout and does\n# not depend on paragraph indices.\n# The one structural change (the deleted-note sentence) also relocates the\n# hidden \"_GoBack\" bookmark so it again sits right after the word \"\u05d0\u05e9\u05e4\u05d4\".\n\n$d = $word.ActiveDocument\n\nfunction Replace-Sentence($doc, $oldText, $newText) {\n    $rng = $doc.Content\n    $found = $rng.Find.Execute($oldText)\n    if (-not $found) {\n        throw \"Could not find expected text: $oldText\"\n    }\n    $rng.Text = $newText\n}\n\n# 1) \"...\u05d5\u05d1\u05e6\u05d3 \u05e9\u05de\u05d0\u05dc \u05d4\u05ea\u05d0\u05e8\u05d9\u05da \u05e9\u05d4\u05d5\u05db\u05e0\u05e1\u05d4 \u05d1\u05d5.\" -> \"...\u05d5\u05d1\u05e6\u05d3 \u05e9\u05de\u05d0\u05dc \u05ea\u05d0\u05e8\u05d9\u05da \u05d4\u05d5\u05e1\u05e4\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4.\"\nReplace-Sentence $d \"\u05d1\u05e6\u05d3 \u05d9\u05de\u05d9\u05df \u05de\u05d5\u05e4\u05d9\u05e2\u05d4 \u05db\u05d5\u05ea\u05e8\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4, \u05d5\u05d1\u05e6\u05d3 \u05e9\u05de\u05d0\u05dc \u05d4\u05ea\u05d0\u05e8\u05d9\u05da \u05e9\u05d4\u05d5\u05db\u05e0\u05e1\u05d4 \u05d1\u05d5.\" \"\u05d1\u05e6\u05d3 \u05d9\u05de\u05d9\u05df \u05de\u05d5\u05e4\u05d9\u05e2\u05d4 \u05db\u05d5\u05ea\u05e8\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4, \u05d5\u05d1\u05e6\u05d3 \u05e9\u05de\u05d0\u05dc \u05ea\u05d0\u05e8\u05d9\u05da \u05d4\u05d5\u05e1\u05e4\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4.\"\n\n# 2) Clicking a note -> describe what the child sees.\nReplace-Sentence $d \"\u05d0\u05dd \u05dc\u05d5\u05d7\u05e6\u05d9\u05dd \u05e2\u05dc \u05d4\u05e2\u05e8\u05d4 \u05de\u05e1\u05d5\u05d9\u05de\u05ea \u05d4\u05d9\u05dc\u05d3 \u05d9\u05db\u05d5\u05dc \u05dc\u05e8\u05d0\u05d5\u05ea \u05d0\u05ea \u05d4\u05ea\u05d5\u05db\u05df, \u05e1\u05e4\u05e8 \u05d0\u05dd \u05d4\u05d5\u05e1\u05d9\u05e3 \u05e1\u05e4\u05e8 \u05de\u05e1\u05d5\u05d9\u05d9\u05dd \u05d5\u05d2\u05dd \u05e6\u05d9\u05d5\u05e8. \u05dc\u05de\u05e9\u05dc:\" \"\u05d0\u05dd \u05dc\u05d5\u05d7\u05e6\u05d9\u05dd \u05e2\u05dc \u05d4\u05e2\u05e8\u05d4 \u05de\u05e1\u05d5\u05d9\u05de\u05ea \u05d4\u05d9\u05dc\u05d3 \u05d9\u05db\u05d5\u05dc \u05dc\u05e8\u05d0\u05d5\u05ea \u05d0\u05ea \u05d4\u05ea\u05d5\u05db\u05df, \u05e9\u05dd \u05d4\u05e1\u05e4\u05e8 \u05d0\u05dd \u05d4\u05d5\u05e1\u05d9\u05e3 \u05e1\u05e4\u05e8 \u05de\u05e1\u05d5\u05d9\u05d9\u05dd \u05d5\u05d2\u05dd \u05e6\u05d9\u05d5\u05e8 \u05d0\u05dd \u05d4\u05e2\u05dc\u05d4 \u05e6\u05d9\u05d5\u05e8. \u05dc\u05de\u05e9\u05dc:\"\n\n# 3) \"\u05d4\u05d5\u05e1\u05e3 \u05d4\u05e2\u05e8\u05d4\" button -> opens blank form wording.\nReplace-Sentence $d \"\u05d0\u05dd \u05d4\u05d9\u05dc\u05d3 \u05dc\u05d5\u05d7\u05e5 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \"\"\u05d4\u05d5\u05e1\u05e3 \u05d4\u05e2\u05e8\u05d4\"\", \u05d4\u05d5\u05d0 \u05de\u05e7\u05d1\u05dc \u05e4\u05d5\u05e8\u05dd \u05db\u05d3\u05d9 \u05dc\u05d4\u05db\u05e0\u05d9\u05e1 \u05e4\u05e8\u05d8\u05d9 \u05d4\u05d4\u05e2\u05e8\u05d4 \u05d4\u05d7\u05d3\u05e9\u05d4:\" \"\u05d0\u05dd \u05d4\u05d9\u05dc\u05d3 \u05dc\u05d5\u05d7\u05e5 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \"\"\u05d4\u05d5\u05e1\u05e3 \u05d4\u05e2\u05e8\u05d4\"\", \u05e0\u05e4\u05ea\u05d7 \u05d8\u05d5\u05e4\u05e1 \u05e8\u05d9\u05e7 \u05dc\u05d4\u05d5\u05e1\u05e4\u05ea \u05d4\u05e2\u05e8\u05d4 \u05d4\u05d7\u05d3\u05e9\u05d4:\"\n\n# 4) Book list bullet wording.\nReplace-Sentence $d \"\u05e1\u05e4\u05e8- \u05d1\u05d7\u05d9\u05e8\u05ea \u05e1\u05e4\u05e8 \u05de\u05e1\u05d5\u05d9\u05dd \u05de\u05e8\u05e9\u05d9\u05de\u05ea \u05d4\u05e1\u05e4\u05e8\u05d9 \u05e9\u05dc\u05d5 \u05d0\u05d5 \u05d1\u05dc\u05d9 \u05e1\u05e4\u05e8.\" \"\u05e1\u05e4\u05e8- \u05d1\u05d7\u05d9\u05e8\u05ea \u05e1\u05e4\u05e8 \u05de\u05e1\u05d5\u05d9\u05dd \u05de\u05e8\u05e9\u05d9\u05de\u05ea \u05d4\u05e1\u05e4\u05e8\u05d9\u05dd \u05e9\u05dc \u05d4\u05d9\u05dc\u05d3 \u05d0\u05d5 \u05d1\u05dc\u05d9 \u05e1\u05e4\u05e8.\"\n\n# 5) Pencil/update button wording.\nReplace-Sentence $d \"\u05d1\u05dc\u05d7\u05d9\u05e6\u05d4 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \u05e2\u05e4\u05e8\u05d5\u05df, \u05d4\u05d9\u05dc\u05d3 \u05de\u05e7\u05d1\u05dc \u05e4\u05d5\u05e8\u05dd \u05db\u05d3\u05d9 \u05dc\u05e2\u05d3\u05db\u05df \u05d0\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4:\" \"\u05d1\u05dc\u05d7\u05d9\u05e6\u05d4 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \u05e2\u05e4\u05e8\u05d5\u05df, \u05d4\u05d9\u05dc\u05d3 \u05de\u05e7\u05d1\u05dc \u05d8\u05d5\u05e4\u05e1 \u05e2\u05dc \u05de\u05e0\u05ea \u05dc\u05e2\u05d3\u05db\u05df \u05d0\u05ea \u05d4\u05d4\u05e2\u05e8\u05d4:\"\n\n# 6) Required-fields note.\nReplace-Sentence $d \"\u05d7\u05d9\u05d9\u05d1 \u05e9\u05d4\u05db\u05d5\u05ea\u05e8\u05ea \u05d5\u05d4\u05ea\u05d5\u05db\u05df \u05dc\u05d0 \u05d9\u05d4\u05d9\u05d5 \u05e8\u05d9\u05e7\u05d9\u05dd.\" \"\u05e9\u05d3\u05d4 \u05d4\u05db\u05d5\u05ea\u05e8\u05ea \u05d5\u05e9\u05d3\u05d4 \u05d4\u05ea\u05d5\u05db\u05df \u05dc\u05d0 \u05d9\u05d4\u05d9\u05d5 \u05e8\u05d9\u05e7\u05d9\u05dd.\"\n\n# 7) Delete-note sentence, plus relocating the hidden \"_GoBack\" bookmark so\n# it ends up right after the new word \"\u05d0\u05e9\u05e4\u05d4\" (matching its position in the\n# edited document instead of the trailing empty paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\nReplace-Sentence $d \"\u05d1\u05dc\u05d7\u05d9\u05e6\u05d4 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \u05e4\u05d7 \u05d4\u05d4\u05e2\u05e8\u05d4 \u05e0\u05de\u05d7\u05e7\u05ea.\" \"\u05d1\u05dc\u05d7\u05d9\u05e6\u05d4 \u05e2\u05dc \u05db\u05e4\u05ea\u05d5\u05e8 \u05d0\u05e9\u05e4\u05d4 \u05d4\u05d4\u05e2\u05e8\u05d4 \u05ea\u05d9\u05de\u05d7\u05e7.\"\n\n$trashRng = $d.Content\n$foundWord = $trashRng.Find.Execute(\"\u05d0\u05e9\u05e4\u05d4\")\nif (-not $foundWord) {\n    throw \"Could not find the word \u05d0\u05e9\u05e4\u05d4 after replacement\"\n}\n$bmRange = $d.Range($trashRng.End, $trashRng.End)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
